# "Generate Report for Handoff"
# The localization status report is regenerated: the still-in-progress
# "In Translation" status becomes "Ready for handoff" (it is now ready to
# hand off to localization), and the associated generate/handoff
# timestamps are refreshed to the new run's time. Because the new status
# text is longer than the old one, the Status column is re-autosized
# (wider) on every sheet that shows it.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" -------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Refreshed timestamps for this handoff run -------------------------
$wsOverview.Range("G2").Value = "2016-11-03 20:02:02"
$wsDeDe.Range("H2").Value     = "2016-11-03 20:02:02"
$wsZhCn.Range("H2").Value     = "2016-11-03 20:01:49"

# --- Widen the Status column to fit "Ready for handoff" -----------------
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332
$wsZhCn.Columns.Item(3).ColumnWidth     = 16.333333333333332
$wsDeDe.Columns.Item(3).ColumnWidth     = 16.333333333333332
